$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7641254067420959
$ws.Range("B1").Value = 1.434107184410095
$ws.Range("C1").Value = 5.46373176574707
$ws.Range("D1").Value = 3.170438766479492
$ws.Range("E1").Value = 1.507601022720337
